# Auto-generated edit script applying cached market-price / profit recalculation updates
# to the Aegis_Profits leve-profit tables across all class sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 556.8461
$ws.Range("I18").Value = 519.9167
$ws.Range("K18").Value = 519.9167
$ws.Range("M18").Value = -235.9167
$ws.Range("H28").Value = 603.8484999999999
$ws.Range("J28").Value = 1095.3572
$ws.Range("L28").Value = 1095.3572
$ws.Range("N28").Value = -2065.3572
$ws.Range("H76").Value = 4461.769
$ws.Range("I76").Value = 4116.6665
$ws.Range("J76").Value = 4757.5713
$ws.Range("K76").Value = 4116.6665
$ws.Range("L76").Value = 4757.5713
$ws.Range("M76").Value = -3801.6665
$ws.Range("N76").Value = -5387.5713
$ws.Range("H79").Value = 4461.769
$ws.Range("I79").Value = 4116.6665
$ws.Range("J79").Value = 4757.5713
$ws.Range("K79").Value = 4116.6665
$ws.Range("L79").Value = 4757.5713
$ws.Range("M79").Value = -3024.6665
$ws.Range("N79").Value = -6941.5713
$ws.Range("H127").Value = 14494921
$ws.Range("I127").Value = 449.14285
$ws.Range("J127").Value = 16131394
$ws.Range("K127").Value = 1347.42855
$ws.Range("L127").Value = 48394182
$ws.Range("M127").Value = 3612.57145
$ws.Range("N127").Value = -48404102
$ws.Range("H129").Value = 3061.3914
$ws.Range("I129").Value = 9549.727999999999
$ws.Range("J129").Value = 1022.2
$ws.Range("K129").Value = 28649.184
$ws.Range("L129").Value = 3066.6
$ws.Range("M129").Value = -23649.184
$ws.Range("N129").Value = -13066.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 203956
$ws.Range("I45").Value = 1000000
$ws.Range("K45").Value = 1000000
$ws.Range("M45").Value = -999623
$ws.Range("H110").Value = 143158000
$ws.Range("I110").Value = 167017330
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 167017330
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -167015285
$ws.Range("N110").Value = -6090
$ws.Range("H122").Value = 2489.25
$ws.Range("I122").Value = 3480
$ws.Range("J122").Value = 1498.5
$ws.Range("K122").Value = 10440
$ws.Range("L122").Value = 4495.5
$ws.Range("M122").Value = -7990
$ws.Range("N122").Value = -9395.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1697
$ws.Range("I99").Value = 2440
$ws.Range("J99").Value = 1449.3334
$ws.Range("K99").Value = 2440
$ws.Range("L99").Value = 1449.3334
$ws.Range("M99").Value = -942
$ws.Range("N99").Value = -4445.3334
$ws.Range("H107").Value = 43496650
$ws.Range("I107").Value = 62525972
$ws.Range("J107").Value = 1050.4286
$ws.Range("K107").Value = 62525972
$ws.Range("L107").Value = 1050.4286
$ws.Range("M107").Value = -62524052
$ws.Range("N107").Value = -4890.4286

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 6008
$ws.Range("I33").Value = 6008
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 6008
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -5629
$ws.Range("N33").ClearContents()
$ws.Range("H122").Value = 689.2
$ws.Range("I122").Value = 373
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 1119
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = 1331
$ws.Range("N122").Value = -7600

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7183.3335
$ws.Range("I87").Value = 4325
$ws.Range("J87").Value = 12900
$ws.Range("K87").Value = 12975
$ws.Range("L87").Value = 38700
$ws.Range("M87").Value = -11727
$ws.Range("N87").Value = -41196
$ws.Range("H90").Value = 7183.3335
$ws.Range("I90").Value = 4325
$ws.Range("J90").Value = 12900
$ws.Range("K90").Value = 38925
$ws.Range("L90").Value = 116100
$ws.Range("M90").Value = -32685
$ws.Range("N90").Value = -128580
$ws.Range("H114").Value = 418.91666
$ws.Range("I114").Value = 181.45454
$ws.Range("J114").Value = 3031
$ws.Range("K114").Value = 544.3636200000001
$ws.Range("L114").Value = 9093
$ws.Range("M114").Value = 2709.63638
$ws.Range("N114").Value = -15601

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6096666.5
$ws.Range("I11").Value = 6644545.5
$ws.Range("J11").Value = 70000
$ws.Range("K11").Value = 6644545.5
$ws.Range("L11").Value = 70000
$ws.Range("M11").Value = -6644406.5
$ws.Range("N11").Value = -70278
$ws.Range("H21").Value = 11600
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9827
$ws.Range("H30").Value = 11600
$ws.Range("I30").Value = 10000
$ws.Range("K30").Value = 10000
$ws.Range("M30").Value = -9895
$ws.Range("H80").Value = 200002020
$ws.Range("I80").Value = 250001970
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 250001970
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -250000972
$ws.Range("N80").Value = -4196
$ws.Range("H83").Value = 200002020
$ws.Range("I83").Value = 250001970
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 1250009850
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -1250004858
$ws.Range("N83").Value = -20984
$ws.Range("H102").Value = 430647
$ws.Range("I102").Value = 2144.889
$ws.Range("J102").Value = 1201950.8
$ws.Range("K102").Value = 2144.889
$ws.Range("L102").Value = 1201950.8
$ws.Range("M102").Value = -522.8890000000001
$ws.Range("N102").Value = -1205194.8
$ws.Range("H113").Value = 1995.591
$ws.Range("I113").Value = 1912.75
$ws.Range("J113").Value = 2042.9286
$ws.Range("K113").Value = 1912.75
$ws.Range("L113").Value = 2042.9286
$ws.Range("M113").Value = 257.25
$ws.Range("N113").Value = -6382.9286
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -10900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 85558.914
$ws.Range("I40").Value = 334800
$ws.Range("J40").Value = 2478.5557
$ws.Range("K40").Value = 334800
$ws.Range("L40").Value = 2478.5557
$ws.Range("M40").Value = -334664
$ws.Range("N40").Value = -2750.5557
$ws.Range("H69").Value = 34949.5
$ws.Range("J69").Value = 34949.5
$ws.Range("L69").Value = 34949.5
$ws.Range("N69").Value = -36571.5
$ws.Range("H72").Value = 34949.5
$ws.Range("J72").Value = 34949.5
$ws.Range("L72").Value = 104848.5
$ws.Range("N72").Value = -112960.5
$ws.Range("H132").Value = 2571.5103
$ws.Range("I132").Value = 2657.1843
$ws.Range("K132").Value = 7971.5529
$ws.Range("M132").Value = -5441.5529

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 90910300
$ws.Range("I96").Value = 200001470
$ws.Range("J96").Value = 1000.1667
$ws.Range("K96").Value = 200001470
$ws.Range("L96").Value = 1000.1667
$ws.Range("M96").Value = -200000097
$ws.Range("N96").Value = -3746.1667
$ws.Range("H122").Value = 4003.3333
$ws.Range("J122").Value = 4003
$ws.Range("L122").Value = 12009
$ws.Range("N122").Value = -16909
$ws.Range("H136").Value = 1052.75
$ws.Range("I136").Value = 670.1818
$ws.Range("J136").Value = 1300.2941
$ws.Range("K136").Value = 2010.5454
$ws.Range("L136").Value = 3900.8823
$ws.Range("M136").Value = 539.4546
$ws.Range("N136").Value = -9000.882300000001
